# Generate Report for Archive
#
# The localization-status workbook tracks each handed-off file's status
# per target language. Two files (2ebdb7a0-3508-44a4-b7b4-dbb8d8ae1e81.md
# and a9e74d13-fdff-4165-ab41-4bcd77bf9071.md) have moved out of the
# "Ready for handoff" state and are now "In Translation" for both the
# zh-cn and de-de locales. Update the Status column (column C) on each
# locale sheet's table accordingly.

$wb = $excel.ActiveWorkbook

$targetFiles = @(
    "2ebdb7a0-3508-44a4-b7b4-dbb8d8ae1e81.md",
    "a9e74d13-fdff-4165-ab41-4bcd77bf9071.md"
)

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column A = "Source File Name", Column C = "Status" (row 1 is the header).
    $row = 2
    while ($ws.Cells.Item($row, 1).Value2 -ne $null) {
        $fileName = $ws.Cells.Item($row, 1).Value2
        if ($targetFiles -contains $fileName) {
            $ws.Cells.Item($row, 3).Value = "In Translation"
        }
        $row = $row + 1
    }
}
